# "Generate Report for Handback" - refresh the generated/handoff/handback
# timestamps that get stamped into the handback-status report each run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$overview.Range("G2").Value = "2016-08-21 03:07:27"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-08-21 03:07:22"
$zhcn.Range("K2").Value = "2016-08-21 03:07:39"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$dede.Range("H2").Value = "2016-08-21 03:07:27"
$dede.Range("K2").Value = "2016-08-21 03:07:46"
